$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: fill in the new model result row ---

# Plain numeric columns.
$ws.Range("A19").Value = 12
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 60
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 5

# Columns are populated left-to-right (F..K) so new shared-string entries
# are appended in the same order the reference workbook used.

# Text columns whose content looks like a plain number / a percentage.
# Force text interpretation first so Excel stores the literal string
# instead of converting it to a numeric value.
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "2.9653518199920654"

# Text column whose content does not resemble a number, so Excel will not
# silently re-type it.
$ws.Range("G19").Value = "3m 8s / 2m 54s"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "12.526968819089731 %"

$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = "5.8480646276361306 %"

$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "16.88605324018926 %"

$ws.Range("K19").Value = "31/12 epochs"

# Re-apply the same formatting as the row above across the whole row so
# every cell keeps the table's normal style (rather than a freshly minted
# "Text" style only for the cells touched above).
$ws.Range("A18:K18").Copy()
$ws.Range("A19:K19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
